# DataEngine.xlsx: fill in the previously-blank "expected result" cells
# with the new "FAIL" status (a sibling value to the existing "PASS"
# shared string) on both sheets.

$wb = $excel.ActiveWorkbook

$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsCases = $wb.Worksheets.Item("Test Cases")

# Test Steps!G12 was empty -> now holds "FAIL"
$wsSteps.Range("G12").Value = "FAIL"

# Test Cases!D3 was empty -> now holds "FAIL"
$wsCases.Range("D3").Value = "FAIL"
